$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.717901
$ws.Range("H2").Value = 2.153703
$ws.Range("I2").Value = 0.0380297505351077
$ws.Range("J2").Value = 0.0380297505351077
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 104.794801
$ws.Range("N2").Value = 314.384403
$ws.Range("O2").Value = 0.3872421191355361
$ws.Range("P2").Value = 0.3872421191355361
$ws.Range("Q2").Value = 75.232292432701
$ws.Range("R2").Value = 677.0906318943091
$ws.Range("S2").Value = 0.0147267211874109
$ws.Range("T2").Value = 0.0147267211874109
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.717901
$ws.Range("H3").Value = 2.153703
$ws.Range("I3").Value = 0.0380297505351077
$ws.Range("J3").Value = 0.0380297505351077
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 136.674446
$ws.Range("N3").Value = 410.023338
$ws.Range("O3").Value = 0.5050451128841343
$ws.Range("P3").Value = 0.5050451128841343
$ws.Range("Q3").Value = 98.11872145784599
$ws.Range("R3").Value = 883.068493120614
$ws.Range("S3").Value = 0.01920673965195894
$ws.Range("T3").Value = 0.01920673965195894
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.717901
$ws.Range("H4").Value = 2.153703
$ws.Range("I4").Value = 0.0380297505351077
$ws.Range("J4").Value = 0.0380297505351077
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 29.14904533333333
$ws.Range("N4").Value = 87.447136
$ws.Range("O4").Value = 0.1077127679803296
$ws.Range("P4").Value = 0.1077127679803296
$ws.Range("Q4").Value = 20.92612879384533
$ws.Range("R4").Value = 188.335159144608
$ws.Range("S4").Value = 0.00409628969573787
$ws.Range("T4").Value = 0.00409628969573787
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 13.91986866666667
$ws.Range("H5").Value = 41.75960600000001
$ws.Range("I5").Value = 0.7373845876726675
$ws.Range("J5").Value = 0.7373845876726675
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 104.794801
$ws.Range("N5").Value = 314.384403
$ws.Range("O5").Value = 0.3872421191355361
$ws.Range("P5").Value = 0.3872421191355361
$ws.Range("Q5").Value = 1458.729866869469
$ws.Range("R5").Value = 13128.56880182522
$ws.Range("S5").Value = 0.2855463703482473
$ws.Range("T5").Value = 0.2855463703482473
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 13.91986866666667
$ws.Range("H6").Value = 41.75960600000001
$ws.Range("I6").Value = 0.7373845876726675
$ws.Range("J6").Value = 0.7373845876726675
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 136.674446
$ws.Range("N6").Value = 410.023338
$ws.Range("O6").Value = 0.5050451128841343
$ws.Range("P6").Value = 0.5050451128841343
$ws.Range("Q6").Value = 1902.490338409425
$ws.Range("R6").Value = 17122.41304568483
$ws.Range("S6").Value = 0.3724124823201632
$ws.Range("T6").Value = 0.3724124823201632
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 13.91986866666667
$ws.Range("H7").Value = 41.75960600000001
$ws.Range("I7").Value = 0.7373845876726675
$ws.Range("J7").Value = 0.7373845876726675
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 29.14904533333333
$ws.Range("N7").Value = 87.447136
$ws.Range("O7").Value = 0.1077127679803296
$ws.Range("P7").Value = 0.1077127679803296
$ws.Range("Q7").Value = 405.750882798713
$ws.Range("R7").Value = 3651.757945188417
$ws.Range("S7").Value = 0.07942573500425702
$ws.Range("T7").Value = 0.07942573500425702
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.239582666666666
$ws.Range("H8").Value = 12.718748
$ws.Range("I8").Value = 0.2245856617922248
$ws.Range("J8").Value = 0.2245856617922248
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 104.794801
$ws.Range("N8").Value = 314.384403
$ws.Range("O8").Value = 0.3872421191355361
$ws.Range("P8").Value = 0.3872421191355361
$ws.Range("Q8").Value = 444.2862218763826
$ws.Range("R8").Value = 3998.575996887444
$ws.Range("S8").Value = 0.08696902759987794
$ws.Range("T8").Value = 0.08696902759987796
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.239582666666666
$ws.Range("H9").Value = 12.718748
$ws.Range("I9").Value = 0.2245856617922248
$ws.Range("J9").Value = 0.2245856617922248
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 136.674446
$ws.Range("N9").Value = 410.023338
$ws.Range("O9").Value = 0.5050451128841343
$ws.Range("P9").Value = 0.5050451128841343
$ws.Range("Q9").Value = 579.4426122378692
$ws.Range("R9").Value = 5214.983510140823
$ws.Range("S9").Value = 0.1134258909120122
$ws.Range("T9").Value = 0.1134258909120122
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.239582666666666
$ws.Range("H10").Value = 12.718748
$ws.Range("I10").Value = 0.2245856617922248
$ws.Range("J10").Value = 0.2245856617922248
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 29.14904533333333
$ws.Range("N10").Value = 87.447136
$ws.Range("O10").Value = 0.1077127679803296
$ws.Range("P10").Value = 0.1077127679803296
$ws.Range("Q10").Value = 123.5797873450809
$ws.Range("R10").Value = 1112.218086105728
$ws.Range("S10").Value = 0.02419074328033468
$ws.Range("T10").Value = 0.02419074328033468
